$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $val) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $val
    $cellRange.ClearFormats()
}

Set-TextValue $ws.Range("D2") '61.227.25'
Set-TextValue $ws.Range("E2") '  -4.09%  '

Set-TextValue $ws.Range("D3") '2.996.02'
Set-TextValue $ws.Range("E3") '  -3.21%  '

Set-TextValue $ws.Range("E4") '  +0.03%  '

Set-TextValue $ws.Range("D5") '536.99'
Set-TextValue $ws.Range("E5") '  -0.94%  '

Set-TextValue $ws.Range("D6") '135.76'

Set-TextValue $ws.Range("D8") '2.990.86'
Set-TextValue $ws.Range("E8") '  -3.11%  '

Set-TextValue $ws.Range("E9") '  -0.43%  '

Set-TextValue $ws.Range("E10") '  -4.93%  '

Set-TextValue $ws.Range("E11") '  +0.31%  '

Set-TextValue $ws.Range("D12") '0.449'
Set-TextValue $ws.Range("E12") '  -2.37%  '

Set-TextValue $ws.Range("D13") '0.0000221'
Set-TextValue $ws.Range("E13") '  -2.60%  '

Set-TextValue $ws.Range("D14") '34.17'
Set-TextValue $ws.Range("E14") '  -1.85%  '

Set-TextValue $ws.Range("D15") '3.488.31'
Set-TextValue $ws.Range("E15") '  -2.99%  '

Set-TextValue $ws.Range("E16") '  -1.27%  '

Set-TextValue $ws.Range("D17") '61.318.91'
Set-TextValue $ws.Range("E17") '  -3.89%  '

Set-TextValue $ws.Range("D18") '3.002.64'
Set-TextValue $ws.Range("E18") '  -3.02%  '

Set-TextValue $ws.Range("D19") '6.64'
Set-TextValue $ws.Range("E19") '  -1.08%  '

Set-TextValue $ws.Range("D20") '467.24'
Set-TextValue $ws.Range("E20") '  -4.73%  '

Set-TextValue $ws.Range("D21") '13.28'
Set-TextValue $ws.Range("E21") '  -1.73%  '

Set-TextValue $ws.Range("D22") '0.676'
Set-TextValue $ws.Range("E22") '  -3.89%  '

Set-TextValue $ws.Range("D23") '6.95'
Set-TextValue $ws.Range("E23") '  -3.72%  '

Set-TextValue $ws.Range("D24") '79.94'
Set-TextValue $ws.Range("E24") '  +0.01%  '

Set-TextValue $ws.Range("D25") '12.03'
Set-TextValue $ws.Range("E25") '  -2.07%  '

Set-TextValue $ws.Range("D26") '0.999'
Set-TextValue $ws.Range("E26") '  -0.34%  '

Set-TextValue $ws.Range("E27") '  -1.97%  '

Set-TextValue $ws.Range("D28") '7.83'
Set-TextValue $ws.Range("E28") '  -6.94%  '

Set-TextValue $ws.Range("E29") '  +0.08%  '

Set-TextValue $ws.Range("D30") '1.89'
Set-TextValue $ws.Range("E30") '  -1.25%  '

Set-TextValue $ws.Range("D31") '1.15'
Set-TextValue $ws.Range("E31") '  +2.94%  '

Set-TextValue $ws.Range("D32") '25.63'
Set-TextValue $ws.Range("E32") '  -2.72%  '

Set-TextValue $ws.Range("D33") '5.50'
Set-TextValue $ws.Range("E33") '  +1.39%  '

Set-TextValue $ws.Range("D34") '55.55'
Set-TextValue $ws.Range("E34") '  -3.29%  '

Set-TextValue $ws.Range("D35") '2.28'
Set-TextValue $ws.Range("E35") '  -5.74%  '

Set-TextValue $ws.Range("D36") '5.90'
Set-TextValue $ws.Range("E36") '  -3.04%  '

Set-TextValue $ws.Range("D37") '454.16'
Set-TextValue $ws.Range("E37") '  -8.71%  '

Set-TextValue $ws.Range("D38") '3.177.55'
Set-TextValue $ws.Range("E38") '  -1.07%  '

Set-TextValue $ws.Range("D39") '0.0790'
Set-TextValue $ws.Range("E39") '  -1.48%  '

Set-TextValue $ws.Range("D40") '0.0386'
Set-TextValue $ws.Range("E40") '  -3.79%  '

Set-TextValue $ws.Range("D41") '0.118'
Set-TextValue $ws.Range("E41") '  -0.25%  '

Set-TextValue $ws.Range("D42") '8.15'
Set-TextValue $ws.Range("E42") '  -0.34%  '

Set-TextValue $ws.Range("D43") '2.49'
Set-TextValue $ws.Range("E43") '  -8.12%  '

Set-TextValue $ws.Range("D44") '27.35'
Set-TextValue $ws.Range("E44") '  +10.31%  '

Set-TextValue $ws.Range("E46") '  -4.78%  '

Set-TextValue $ws.Range("D47") '2.00'
Set-TextValue $ws.Range("E47") '  -3.07%  '

Set-TextValue $ws.Range("D48") '119.70'
Set-TextValue $ws.Range("E48") '  -1.45%  '

Set-TextValue $ws.Range("E49") '  -1.07%  '

Set-TextValue $ws.Range("D50") '0.0₃0497'
Set-TextValue $ws.Range("E50") '  -9.07%  '

Set-TextValue $ws.Range("B51") 'BitgetToken'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
Set-TextValue $ws.Range("D51") '1.24'
Set-TextValue $ws.Range("E51") '  +5.54%  '
